# Adding csv version of site data
# The edit re-numbers the "site_num" column (A) starting at row 6: a new
# site row (row 6) is inserted into the numbering sequence, so A6 becomes 5
# and every subsequent site_num (rows 7-34, previously 5-32) is bumped up
# by one (6-33). The cells/styles in the other columns are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 previously had no value in column A (site_num) - fill it in with 5,
# then shift every following site_num down the list by incrementing it by 1
# (row 7: 5->6, row 8: 6->7, ... row 34: 32->33).
for ($r = 6; $r -le 34; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Reflect the author's final cursor position/selection in the sheet.
$ws.Range("G13").Select()
